$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14
$ws.Range("A14").Value = 111908768
$ws.Range("B14").Value = 96735
$ws.Range("D14").Value = "VU"
$ws.Range("E14").Value = 220787
$ws.Range("F14").Value = "Knärot"
$ws.Range("G14").Value = "Goodyera repens"
$ws.Range("H14").Value = "(L.) R. Br."
$ws.Range("I14").Value = "'1"
$ws.Range("Q14").Value = 467912
$ws.Range("R14").Value = 6875299

# Row 15
$ws.Range("A15").Value = 111909766
$ws.Range("B15").Value = 89331
$ws.Range("D15").Value = "LC"
$ws.Range("E15").Value = 3215
$ws.Range("F15").Value = "Rödgul trumpetsvamp"
$ws.Range("G15").Value = "Craterellus lutescens"
$ws.Range("H15").Value = "(Fr.) Fr."
$ws.Range("P15").Value = "Fläcksberget, Hjd"
$ws.Range("Q15").Value = 467757
$ws.Range("R15").Value = 6875470

# Row 16
$ws.Range("A16").Value = 111908364
$ws.Range("B16").Value = 90808
$ws.Range("E16").Value = 4362
$ws.Range("F16").Value = "Blå taggsvamp"
$ws.Range("G16").Value = "Hydnellum caeruleum"
$ws.Range("H16").Value = "(Hornem.) P.Karst."
$ws.Range("P16").Value = "Gröbäcken, Hjd"
$ws.Range("Q16").Value = 467724
$ws.Range("R16").Value = 6874811

# Row 17
$ws.Range("A17").Value = 111908700
$ws.Range("B17").Value = 90113
$ws.Range("E17").Value = 760
$ws.Range("F17").Value = "Doftticka"
$ws.Range("G17").Value = "Haploporus odorus"
$ws.Range("H17").Value = "(Sommerf.) Bondartsev & Singer"
$ws.Range("I17").Value = "'6"
$ws.Range("J17").Value = "fruktkroppar"
$ws.Range("Q17").Value = 467922
$ws.Range("R17").Value = 6875307
$ws.Range("AC17").Value = "Förekomst av doftticka i avverkningsanmält område."

# Row 18
$ws.Range("A18").Value = 111909536
$ws.Range("B18").Value = 77402
$ws.Range("D18").Value = "NT"
$ws.Range("E18").Value = 6446
$ws.Range("F18").Value = "Kolflarnlav"
$ws.Range("G18").Value = "Carbonicola anthracophila"
$ws.Range("H18").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("I18").ClearContents()
$ws.Range("J18").ClearContents()
$ws.Range("Q18").Value = 467891
$ws.Range("R18").Value = 6875425
$ws.Range("AC18").ClearContents()

# Row 19
$ws.Range("A19").Value = 111909174
$ws.Range("B19").Value = 77402
$ws.Range("D19").Value = "NT"
$ws.Range("E19").Value = 6446
$ws.Range("F19").Value = "Kolflarnlav"
$ws.Range("G19").Value = "Carbonicola anthracophila"
$ws.Range("H19").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q19").Value = 467989
$ws.Range("R19").Value = 6875353

# Row 20
$ws.Range("B20").Value = 77402

# Row 21
$ws.Range("B21").Value = 90826

# Row 22
$ws.Range("A22").Value = 112014142
$ws.Range("B22").Value = 90814
$ws.Range("D22").Value = "LC"
$ws.Range("E22").Value = 4364
$ws.Range("F22").Value = "Dropptaggsvamp"
$ws.Range("G22").Value = "Hydnellum ferrugineum"
$ws.Range("H22").Value = "(Fr.:Fr.) P. Karst."
$ws.Range("Q22").Value = 467443
$ws.Range("R22").Value = 6875337

# Row 23
$ws.Range("A23").Value = 112014229
$ws.Range("B23").Value = 90830
$ws.Range("E23").Value = 2059
$ws.Range("F23").Value = "Skrovlig taggsvamp"
$ws.Range("G23").Value = "Hydnellum scabrosum"
$ws.Range("H23").Value = "(Fr.) E.Larss., K.H.Larss. & Kõljalg"
$ws.Range("Q23").Value = 467427
$ws.Range("R23").Value = 6875290

# Row 24
$ws.Range("A24").Value = 112014300
$ws.Range("B24").Value = 90837
$ws.Range("Q24").Value = 467415
$ws.Range("R24").Value = 6875287

# Row 25
$ws.Range("A25").Value = 112014923
$ws.Range("B25").Value = 90837
$ws.Range("E25").Value = 5966
$ws.Range("F25").Value = "Motaggsvamp"
$ws.Range("G25").Value = "Sarcodon squamosus"
$ws.Range("H25").Value = "(Schaeff.) Quél."
$ws.Range("Q25").Value = 467413
$ws.Range("R25").Value = 6875234

# Row 26
$ws.Range("A26").Value = 112014208
$ws.Range("B26").Value = 90806
$ws.Range("Q26").Value = 467418
$ws.Range("R26").Value = 6875313

# Row 27
$ws.Range("A27").Value = 112015011
$ws.Range("B27").Value = 90806
$ws.Range("E27").Value = 4361
$ws.Range("F27").Value = "Orange taggsvamp"
$ws.Range("G27").Value = "Hydnellum aurantiacum"
$ws.Range("H27").Value = "(Batsch:Fr.) P.Karst."
$ws.Range("Q27").Value = 467390
$ws.Range("R27").Value = 6875328

# Row 28
$ws.Range("A28").Value = 112014423
$ws.Range("B28").Value = 90806
$ws.Range("Q28").Value = 467430
$ws.Range("R28").Value = 6875238

# Row 29
$ws.Range("A29").Value = 112014177
$ws.Range("B29").Value = 90837
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 5966
$ws.Range("F29").Value = "Motaggsvamp"
$ws.Range("G29").Value = "Sarcodon squamosus"
$ws.Range("H29").Value = "(Schaeff.) Quél."
$ws.Range("Q29").Value = 467390
$ws.Range("R29").Value = 6875328
